$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-12-31 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-01-01 Monday", 2) | Out-Null
$d.Content.Find.Execute("71÷4=17, 3", $true, $false, $false, $false, $false, $true, 1, $false, "69÷7=9, 6", 2) | Out-Null
$d.Content.Find.Execute("19÷8=2, 3", $true, $false, $false, $false, $false, $true, 1, $false, "39÷4=9, 3", 2) | Out-Null
$d.Content.Find.Execute("36÷6=6, 0", $true, $false, $false, $false, $false, $true, 1, $false, "13÷4=3, 1", 2) | Out-Null
$d.Content.Find.Execute("29÷6=4, 5", $true, $false, $false, $false, $false, $true, 1, $false, "25÷3=8, 1", 2) | Out-Null
$d.Content.Find.Execute("80÷8=10, 0", $true, $false, $false, $false, $false, $true, 1, $false, "77÷9=8, 5", 2) | Out-Null
$d.Content.Find.Execute("63÷9=7, 0", $true, $false, $false, $false, $false, $true, 1, $false, "54÷2=27, 0", 2) | Out-Null
$d.Content.Find.Execute("68÷5=13, 3", $true, $false, $false, $false, $false, $true, 1, $false, "93÷2=46, 1", 2) | Out-Null
$d.Content.Find.Execute("20÷8=2, 4", $true, $false, $false, $false, $false, $true, 1, $false, "44÷4=11, 0", 2) | Out-Null
$d.Content.Find.Execute("37÷8=4, 5", $true, $false, $false, $false, $false, $true, 1, $false, "50÷3=16, 2", 2) | Out-Null
$d.Content.Find.Execute("99÷9=11, 0", $true, $false, $false, $false, $false, $true, 1, $false, "77÷4=19, 1", 2) | Out-Null
$d.Content.Find.Execute("69÷3=23, 0", $true, $false, $false, $false, $false, $true, 1, $false, "16÷9=1, 7", 2) | Out-Null
$d.Content.Find.Execute("42÷4=10, 2", $true, $false, $false, $false, $false, $true, 1, $false, "33÷2=16, 1", 2) | Out-Null
$d.Content.Find.Execute("30÷4=7, 2", $true, $false, $false, $false, $false, $true, 1, $false, "98÷5=19, 3", 2) | Out-Null
$d.Content.Find.Execute("87÷8=10, 7", $true, $false, $false, $false, $false, $true, 1, $false, "62÷3=20, 2", 2) | Out-Null
$d.Content.Find.Execute("89÷5=17, 4", $true, $false, $false, $false, $false, $true, 1, $false, "62÷3=20, 2", 2) | Out-Null
$d.Content.Find.Execute("91÷7=13, 0", $true, $false, $false, $false, $false, $true, 1, $false, "77÷6=12, 5", 2) | Out-Null
$d.Content.Find.Execute("88÷7=12, 4", $true, $false, $false, $false, $false, $true, 1, $false, "37÷8=4, 5", 2) | Out-Null
$d.Content.Find.Execute("61÷2=30, 1", $true, $false, $false, $false, $false, $true, 1, $false, "51÷3=17, 0", 2) | Out-Null
$d.Content.Find.Execute("34÷2=17, 0", $true, $false, $false, $false, $false, $true, 1, $false, "83÷3=27, 2", 2) | Out-Null
$d.Content.Find.Execute("79÷5=15, 4", $true, $false, $false, $false, $false, $true, 1, $false, "44÷9=4, 8", 2) | Out-Null
$d.Content.Find.Execute("86÷6=14, 2", $true, $false, $false, $false, $false, $true, 1, $false, "79÷9=8, 7", 2) | Out-Null
$d.Content.Find.Execute("86÷3=28, 2", $true, $false, $false, $false, $false, $true, 1, $false, "96÷6=16, 0", 2) | Out-Null
$d.Content.Find.Execute("59÷3=19, 2", $true, $false, $false, $false, $false, $true, 1, $false, "20÷2=10, 0", 2) | Out-Null
$d.Content.Find.Execute("14÷6=2, 2", $true, $false, $false, $false, $false, $true, 1, $false, "35÷8=4, 3", 2) | Out-Null
$d.Content.Find.Execute("92÷6=15, 2", $true, $false, $false, $false, $false, $true, 1, $false, "83÷4=20, 3", 2) | Out-Null

Write-Host "Replacements complete"